# Round ConvexHullArea values (column D, rows 2-26) to the nearest integer.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $current = $cell.Value2
    if ($current -ne $null) {
        $d = [double]$current
        if ($d -ge 0) {
            $rounded = [Math]::Floor($d + 0.5)
        } else {
            $rounded = [Math]::Ceiling($d - 0.5)
        }
        $cell.Value2 = $rounded
    }
}
